# Update cryptocurrency symbol list prices/volumes (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","E18","D19","E19","E20","D21","E21","D22","E22","D23","E23","E24","D25","E25","D26","E26","D27","E27","D39","E39","D40","E40","D41","E41","D42","E42","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $cells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "304.56"
$ws.Range("E2").Value = "1.00%"
$ws.Range("D3").Value = "35.73"
$ws.Range("E3").Value = "1.67%"
$ws.Range("D4").Value = "5.059"
$ws.Range("E4").Value = "0.42%"
$ws.Range("D5").Value = "0.08044"
$ws.Range("E5").Value = "0.75%"
$ws.Range("D6").Value = "1.924"
$ws.Range("E6").Value = "1.58%"
$ws.Range("D8").Value = "7.733"
$ws.Range("E8").Value = "-0.87%"
$ws.Range("D9").Value = "0.9300"
$ws.Range("E9").Value = "0.80%"
$ws.Range("D10").Value = "0.1386"
$ws.Range("E10").Value = "9.63%"
$ws.Range("D11").Value = "0.1899"
$ws.Range("E11").Value = "2.50%"
$ws.Range("D12").Value = "0.09152"
$ws.Range("E12").Value = "-8.74%"
$ws.Range("D13").Value = "0.03629"
$ws.Range("E13").Value = "2.86%"
$ws.Range("D14").Value = "0.09802"
$ws.Range("E14").Value = "-0.43%"
$ws.Range("D15").Value = "0.001418"
$ws.Range("E15").Value = "1.95%"
$ws.Range("D16").Value = "0.005903"
$ws.Range("E16").Value = "0.24%"
$ws.Range("D17").Value = "3.555"
$ws.Range("E17").Value = "1.44%"
$ws.Range("E18").Value = "1.20%"
$ws.Range("D19").Value = "0.3468"
$ws.Range("E19").Value = "1.99%"
$ws.Range("E20").Value = "2.36%"
$ws.Range("D21").Value = "4.896"
$ws.Range("E21").Value = "-2.74%"
$ws.Range("D22").Value = "0.2508"
$ws.Range("E22").Value = "4.50%"
$ws.Range("D23").Value = "0.04443"
$ws.Range("E23").Value = "-1.19%"
$ws.Range("E24").Value = "0.70%"
$ws.Range("D25").Value = "0.004785"
$ws.Range("E25").Value = "-0.05%"
$ws.Range("D26").Value = "0.0001561"
$ws.Range("E26").Value = "24.75%"
$ws.Range("D27").Value = "0.0003130"
$ws.Range("E27").Value = "4.31%"
$ws.Range("D39").Value = "0.01960"
$ws.Range("E39").Value = "4.03%"
$ws.Range("D40").Value = "0.04892"
$ws.Range("E40").Value = "3.63%"
$ws.Range("D41").Value = "0.007624"
$ws.Range("E41").Value = "1.52%"
$ws.Range("D42").Value = "0.009194"
$ws.Range("E42").Value = "-9.85%"
$ws.Range("E43").Value = "3.69%"
$ws.Range("D44").Value = "0.002101"
$ws.Range("E44").Value = "-0.51%"
$ws.Range("D45").Value = "0.01134"
$ws.Range("E45").Value = "7.39%"
$ws.Range("D46").Value = "0.00006389"
$ws.Range("E46").Value = "2.70%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "0.04%"
$ws.Range("D48").Value = "63.57"
$ws.Range("E48").Value = "-1.41%"
$ws.Range("D49").Value = "0.001191"
$ws.Range("E49").Value = "-20.03%"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "0.04%"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "0.04%"
